$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "S 11" week row
$ws.Range("D15").Value = "S 11"
$ws.Range("E15").Value = 44374
$ws.Range("F15").Value = 44381

# New "S 12" week row
$ws.Range("D16").Value = "S 12"
$ws.Range("E16").Value = 44388
$ws.Range("F16").Value = 44395

# E15 / F15 / F16 reuse the same date style already used by the rest of the
# column (copy the formatting instead of re-deriving it, so the existing
# style index is reused rather than a new duplicate one being minted).
$ws.Range("E5:F5").Copy()
$ws.Range("E15:F15").PasteSpecial(-4122)
$ws.Range("F6").Copy()
$ws.Range("F16").PasteSpecial(-4122)

# E16 gets a distinct day-month display (the "real time" cell) with the same
# left/center alignment as the other date cells.
$ws.Range("E16").NumberFormat = "d-mmm"
$ws.Range("E16").HorizontalAlignment = -4131
$ws.Range("E16").VerticalAlignment = -4108

# Selection moves to E15 after the edit
$null = $ws.Range("E15").Select()

# Printing layout settings added alongside the new rows
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
